# Applies the "ifo GDP component analysis preprocessing" update:
# extends the staircase-shaped error matrix on Sheet1 by one more
# anti-diagonal of values, and updates the previously-last anti-diagonal
# (row+col = 17) to the newly recomputed averages (row+col = 18 is new).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cells (previous trailing diagonal, now recomputed) ---
$ws.Range("K6").Value  = 1.388747888886706
$ws.Range("J7").Value  = 0.444773652920949
$ws.Range("I8").Value  = 0.4487415504340581
$ws.Range("H9").Value  = 0.5843816406042994
$ws.Range("G10").Value = 0.3435754587486348
$ws.Range("F11").Value = 0.2982442434965384
$ws.Range("E12").Value = 0.2313828215604846
$ws.Range("D13").Value = 0.201796619203768
$ws.Range("C14").Value = 0.1836459624741271
$ws.Range("B15").Value = 0.1656141382254278

# --- Add newly populated cells (new trailing diagonal) ---
$ws.Range("K7").Value  = 0.2348700177716323
$ws.Range("J8").Value  = 0.2388379152847414
$ws.Range("I9").Value  = 0.3744780054549828
$ws.Range("H10").Value = 0.1336718235993181
$ws.Range("G11").Value = 0.08834060834722172
$ws.Range("F12").Value = 0.02147918641116785
$ws.Range("E13").Value = -0.00810701594554874
$ws.Range("D14").Value = -0.02625767267518964
$ws.Range("C15").Value = -0.04428949692388896
$ws.Range("B16").Value = -0.09587373626955231
